# Auto update stock data
# Updates the "Date_1" column (A) from 2025/11/16 to 2025/11/17 for the
# first row of each of the 13 companies, and corrects the EBITDA value
# for Alro Steel (row 38) from 38.59 to 32.35.
#
# The date/number-looking values are stored as literal text in this
# workbook, so we prefix each value with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting to a date
# serial number or a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 1).Value = "'2025/11/17"
}

# Alro Steel's EBITDA value changed from 38.59 to 32.35
$ws.Cells.Item(38, 2).Value = "'32.35"
